$p = $ppt.ActivePresentation

# Slide 6 (sldId=261) already has a notes page (notesSlide5.xml) with an
# empty body placeholder -- fill in the speaker notes text.
$s6 = $p.Slides.Item(6)
$np6 = $s6.NotesPage
$np6.Shapes.Item(2).TextFrame.TextRange.Text = "Show how data was acquired from the SRAM and what data was required and why"

# Slide 7 (sldId=262) has no notes page yet -- create the notes body
# placeholder and set its text.
$s7 = $p.Slides.Item(7)
$np7 = $s7.NotesPage
$ph7 = $np7.Shapes.AddPlaceholder(2)
$ph7.TextFrame.TextRange.Text = "Briefly summarise the coding done to calculate hamming distance and carry out the appropriate statistics"

# Slide 8 (sldId=263) has no notes page yet -- create the notes body
# placeholder and set its text.
$s8 = $p.Slides.Item(8)
$np8 = $s8.NotesPage
$ph8 = $np8.Shapes.AddPlaceholder(2)
$ph8.TextFrame.TextRange.Text = "Show figures from the report and discuss what can be inferred from them"

# Slide 9 (sldId=264) has no notes page yet -- create the notes body
# placeholder and set its text.
$s9 = $p.Slides.Item(9)
$np9 = $s9.NotesPage
$ph9 = $np9.Shapes.AddPlaceholder(2)
$ph9.TextFrame.TextRange.Text = "Argue as to whether a microbit is a viable device for creating PUFs"
